$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.297.61"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.823.93"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.48"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4493"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3789"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07440"
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8870"
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.98"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "1.820.92"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.733"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.451"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.48"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07121"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008817"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").Value = "27.310.54"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.390"
$ws.Range("E22").Value = "  +4.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.966"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.67"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.308"
$ws.Range("E26").Value = "  +3.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.64"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.385"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.74"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08905"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7900"
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.202"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.616"
$ws.Range("E33").Value = "  +3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.909"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.112"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01979"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05290"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.366"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5326"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.870"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1713"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.315"
$ws.Range("E43").Value = "  +17.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.652"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5065"
$ws.Range("E45").Value = "  -3.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.67"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.696"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.09"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06398"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.16"
$ws.Range("E51").Value = "  +4.52%  "
